# Add a new "canonical SMILES" column (column D) to the microstate list
# sheet, mirroring "canonical isomeric SMILES" (column C) for every row
# except SM04_micro009, which gets its own distinct canonical SMILES
# (no stereo slashes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D2").Value = "canonical SMILES"

# Column D width (character units), matching the new column in the diff.
$ws.Columns.Item(4).ColumnWidth = 36.85546875

# Row 3 (SM04_micro002)
$ws.Range("D3").Value = "c1ccc2c(c1)c(=[NH+]Cc3ccc(cc3)Cl)nc[nH]2"

# Row 4 (SM04_micro003)
$ws.Range("D4").Value = "c1ccc2c(c1)c(ncn2)NCc3ccc(cc3)Cl"

# Row 5 (SM04_micro004)
$ws.Range("D5").Value = "c1ccc2c(c1)c(=NCc3ccc(cc3)Cl)nc[nH]2"

# Row 6 (SM04_micro005)
$ws.Range("D6").Value = "c1ccc2c(c1)c(ncn2)[N-]Cc3ccc(cc3)Cl"

# Row 7 (SM04_micro006)
$ws.Range("D7").Value = "c1ccc2c(c1)c([nH+]cn2)NCc3ccc(cc3)Cl"

# Row 8 (SM04_micro008)
$ws.Range("D8").Value = "c1ccc2c(c1)c(ncn2)[NH2+]Cc3ccc(cc3)Cl"

# Row 9 (SM04_micro009) - new canonical (non-isomeric) SMILES, different
# from the isomeric SMILES shown in column C.
$ws.Range("D9").Value = "c1ccc2c(c1)c(=NCc3ccc(cc3)Cl)[nH]cn2"

# Row 10 (SM04_micro013)
$ws.Range("D10").Value = "c1ccc2c(c1)c(nc[nH+]2)[NH2+]Cc3ccc(cc3)Cl"

# Row 11 (SM04_micro014)
$ws.Range("D11").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)NCc3ccc(cc3)Cl"

# Row 12 (SM04_micro015)
$ws.Range("D12").Value = "c1ccc2c(c1)c([nH+]cn2)[NH2+]Cc3ccc(cc3)Cl"

# Row 13 (SM04_micro016)
$ws.Range("D13").Value = "c1ccc2c(c1)c([nH+]c[nH+]2)[NH2+]Cc3ccc(cc3)Cl"
